$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new value in C5
$ws.Range("C5").Value = "tester1"

# Update the active cell / selection to C5
$ws.Range("C5").Select()
